# Add the missing "bypassGlobal" (column T) values of -1 to every data row
# that has a Counter (column S) value but no bypassGlobal entry yet, on both
# the arr_ccf and arr_atom sheets. Also restore the view/selection state
# recorded in the saved workbook.

$wb = $excel.ActiveWorkbook

# ---- arr_ccf ----------------------------------------------------------
$wsCcf = $wb.Worksheets.Item("arr_ccf")

$ccfRows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,57,58,59,60,61,62,63,64,68,69,70,71,72,73,76,77,80,81,82,83,84,85,89,90,91,92,93,94,97,98)

foreach ($r in $ccfRows) {
    $wsCcf.Cells.Item($r, 20).Value = -1   # column T = 20
}

# ---- arr_atom -----------------------------------------------------------
$wsAtom = $wb.Worksheets.Item("arr_atom")

$atomRows = @(88,89,90,91,94,95,98,99,100,101,105,106,107,108,111,112)

foreach ($r in $atomRows) {
    $wsAtom.Cells.Item($r, 20).Value = -1   # column T = 20
}

# ---- restore view / selection state -------------------------------------

# arr_atom: scrolled back to the top (E1) with T105 selected
$wsAtom.Activate()
$wsAtom.Range("T105").Select()

# arr_ccf: frozen header row stays frozen; scroll down so row 26 is the
# first visible row under the freeze, with S15 selected
$wsCcf.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$wsCcf.Range("S15").Select()

$wsCcf.Activate()
